# Generate Report for Handoff
# Adds two new localized files (e329ed08-...md and e4b128cf-...md) to the
# localization-status report, pushing the existing ".localization-config"
# row down by two rows on every sheet (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$commitMd    = "c053265932b0b14eeb73e17b5af893534456ddd3"
$commitZh    = "bc65390bf6b2ee9cb3c445a7a5d1093d58dcc930"
$commitDe    = "77b40c8bca0b243664b31f667874eb161e523a3e"

$file1 = "e329ed08-4084-4799-9ae8-3c26ba335479.md"
$file2 = "e4b128cf-6e74-487c-a63b-63c0b070a7ae.md"

$xlfZh1 = "e329ed08-4084-4799-9ae8-3c26ba335479.39805055960258112f342c2ed773a55630a474f6.zh-cn.xlf"
$xlfZh2 = "e4b128cf-6e74-487c-a63b-63c0b070a7ae.bccdea2384581b097eadb01ee47e5801e8184cb5.zh-cn.xlf"
$xlfDe1 = "e329ed08-4084-4799-9ae8-3c26ba335479.39805055960258112f342c2ed773a55630a474f6.de-de.xlf"
$xlfDe2 = "e4b128cf-6e74-487c-a63b-63c0b070a7ae.bccdea2384581b097eadb01ee47e5801e8184cb5.de-de.xlf"

$dtZh = "2016-02-16 07:34:22"
$dtDe = "2016-02-16 07:34:33"
$epoch = "0001-01-01 00:00:00"

$statusReady = "Ready for handoff"
$statusNotLoc = "Not to be localized"

function Set-HyperlinkCell($ws, $cellRef, $text, $url) {
    $ws.Range($cellRef).Value = $text
    $ws.Range($cellRef).Style = "HyperLink"
}

# ---------------------------------------------------------------------
# Sheet "Overview": columns File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 4 used to be ".localization-config" -> becomes the first new file
$ov.Range("B4").Value = $statusReady
$ov.Range("C4").Value = $statusReady

# Row 5 (new): second new file
$ov.Range("B5").Value = $statusReady
$ov.Range("C5").Value = $statusReady

# Row 6 (new): ".localization-config" moves here
$ov.Range("B6").Value = $statusNotLoc
$ov.Range("C6").Value = $statusNotLoc

# Rebuild the hyperlinks for column A (1..3 existing + 3 new/changed ones)
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitMd/e2e/3afcb3a5-4980-43cb-9abd-59c8cdfef388.md", "", "", "3afcb3a5-4980-43cb-9abd-59c8cdfef388.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitMd/e2e/7e85abdc-023e-4001-a7e5-cfc2112e0687.md", "", "", "7e85abdc-023e-4001-a7e5-cfc2112e0687.md")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitMd/e2e/$file1", "", "", $file1)
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitMd/e2e/$file2", "", "", $file2)
$ov.Hyperlinks.Add($ov.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitMd/.localization-config", "", "", ".localization-config")

$ov.Range("A4:A6").Style = "HyperLink"

# ---------------------------------------------------------------------
# Sheets "zh-cn" and "de-de": columns
# A Source File Name | B Status | C Latest Handoff File | D Latest Handoff Datetime |
# E Latest Target File | F Latest Handback File | G Latest Handback DateTime |
# H Handoff Reason | I Dependency From | J Error Detail
# ---------------------------------------------------------------------
function Update-LangSheet($ws, $commitMdHash, $commitXlfHash, $lang, $xlf1, $xlf2, $dt) {
    # Row 4: was ".localization-config" placeholder row, now the 1st new file
    $ws.Range("B4").Value = $statusReady
    $ws.Range("C4").Value = $xlf1
    $ws.Range("D4").Value = $dt
    $ws.Range("G4").Value = $epoch
    $ws.Range("H4").Value = "Include"

    # Row 5 (new): 2nd new file
    $ws.Range("A5").Value = $file2
    $ws.Range("B5").Value = $statusReady
    $ws.Range("C5").Value = $xlf2
    $ws.Range("D5").Value = $dt
    $ws.Range("G5").Value = $epoch
    $ws.Range("H5").Value = "Include"

    # Row 6 (new): ".localization-config" moved here
    $ws.Range("A6").Value = ".localization-config"
    $ws.Range("B6").Value = $statusNotLoc
    $ws.Range("D6").Value = $epoch
    $ws.Range("G6").Value = $epoch
    $ws.Range("H6").Value = "Ignored"

    # Rebuild hyperlinks in column A and C
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitMdHash/e2e/3afcb3a5-4980-43cb-9abd-59c8cdfef388.md", "", "", "3afcb3a5-4980-43cb-9abd-59c8cdfef388.md")
    $ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitXlfHash/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/yuwzho/3afcb3a5-4980-43cb-9abd-59c8cdfef388.bb8e4420bf044dca9fa51faa7f8a0bf0de3c07ee.$lang.xlf", "", "", "3afcb3a5-4980-43cb-9abd-59c8cdfef388.bb8e4420bf044dca9fa51faa7f8a0bf0de3c07ee.$lang.xlf")
    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitMdHash/e2e/7e85abdc-023e-4001-a7e5-cfc2112e0687.md", "", "", "7e85abdc-023e-4001-a7e5-cfc2112e0687.md")
    $ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitXlfHash/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/yuwzho/7e85abdc-023e-4001-a7e5-cfc2112e0687.759644fbffe3e49ac7dffe277ff7e4f735d6b79d.$lang.xlf", "", "", "7e85abdc-023e-4001-a7e5-cfc2112e0687.759644fbffe3e49ac7dffe277ff7e4f735d6b79d.$lang.xlf")
    $ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitMdHash/e2e/$file1", "", "", $file1)
    $ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitXlfHash/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/yuwzho/$xlf1", "", "", $xlf1)
    $ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitMdHash/e2e/$file2", "", "", $file2)
    $ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitXlfHash/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/yuwzho/$xlf2", "", "", $xlf2)
    $ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitMdHash/.localization-config", "", "", ".localization-config")

    $ws.Range("A4:A6").Style = "HyperLink"
    $ws.Range("C4:C5").Style = "HyperLink"
    $ws.Range("D4:D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("G4:G6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

$zh = $wb.Worksheets.Item("zh-cn")
Update-LangSheet $zh $commitMd $commitZh "zh-cn" $xlfZh1 $xlfZh2 $dtZh

$de = $wb.Worksheets.Item("de-de")
Update-LangSheet $de $commitMd $commitDe "de-de" $xlfDe1 $xlfDe2 $dtDe
